$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. All existing data rows (2..461) had their "Förändrad" (column C) date
#    bumped from 2023-09-06 (45175) to 2023-09-08 (45177).
for ($r = 2; $r -le 461; $r++) {
    $ws.Cells.Item($r, 3).Value = 45177
}

# 2. Row 461 picks up an explicit row height (matches the author re-saving
#    the sheet after the edit — Excel stamps ht/customHeight on touched rows).
$ws.Rows(461).RowHeight = 15

# 3. Two brand-new announcements were appended at the bottom of the sheet.

# Row 462 — "A 41460-2023"
$ws.Cells.Item(462, 1).Value = "A 41460-2023"
$ws.Cells.Item(462, 2).Value = 45175
$ws.Cells.Item(462, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(462, 3).Value = 45177
$ws.Cells.Item(462, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(462, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(462, 5).Value = "VÄRNAMO"
$ws.Cells.Item(462, 7).Value = 4
$ws.Cells.Item(462, 8).Value = 0
$ws.Cells.Item(462, 9).Value = 0
$ws.Cells.Item(462, 10).Value = 0
$ws.Cells.Item(462, 11).Value = 0
$ws.Cells.Item(462, 12).Value = 0
$ws.Cells.Item(462, 13).Value = 0
$ws.Cells.Item(462, 14).Value = 0
$ws.Cells.Item(462, 15).Value = 0
$ws.Cells.Item(462, 16).Value = 0
$ws.Cells.Item(462, 17).Value = 0
$ws.Cells.Item(462, 18).WrapText = $true
$ws.Rows(462).RowHeight = 15

# Row 463 — "A 41699-2023"
$ws.Cells.Item(463, 1).Value = "A 41699-2023"
$ws.Cells.Item(463, 2).Value = 45176
$ws.Cells.Item(463, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(463, 3).Value = 45177
$ws.Cells.Item(463, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(463, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(463, 5).Value = "VÄRNAMO"
$ws.Cells.Item(463, 6).Value = "Sveaskog"
$ws.Cells.Item(463, 7).Value = 2.5
$ws.Cells.Item(463, 8).Value = 0
$ws.Cells.Item(463, 9).Value = 0
$ws.Cells.Item(463, 10).Value = 0
$ws.Cells.Item(463, 11).Value = 0
$ws.Cells.Item(463, 12).Value = 0
$ws.Cells.Item(463, 13).Value = 0
$ws.Cells.Item(463, 14).Value = 0
$ws.Cells.Item(463, 15).Value = 0
$ws.Cells.Item(463, 16).Value = 0
$ws.Cells.Item(463, 17).Value = 0
$ws.Cells.Item(463, 18).WrapText = $true
